# "added 4wk low sales check" -- refreshed the forecast numbers (and the
# derived inventory coverage / seasonality / urgency flags that ride along
# with them) on the "Forecast Comparison" sheet, then rolled the new totals
# up into the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Columns: D=MyForecast, H=Inventory Coverage, I=Stockout Risk,
#          J=Reorder Urgency, L=Seasonality Index

# Row 2 (W10)
$wsForecast.Range("D2").Value = 471
$wsForecast.Range("H2").Value = 13.13
$wsForecast.Range("L2").Value = 0.83

# Row 3 (W11)
$wsForecast.Range("D3").Value = 419
$wsForecast.Range("H3").Value = 13.64
$wsForecast.Range("L3").Value = 0.91

# Row 4 (W12)
$wsForecast.Range("D4").Value = 386
$wsForecast.Range("H4").Value = 13.72
$wsForecast.Range("L4").Value = 0.99

# Row 5 (W13)
$wsForecast.Range("D5").Value = 351
$wsForecast.Range("H5").Value = 13.99
$wsForecast.Range("L5").Value = 0.98

# Row 6 (W14)
$wsForecast.Range("D6").Value = 339
$wsForecast.Range("H6").Value = 13.45
$wsForecast.Range("L6").Value = 0.97

# Row 7 (W15)
$wsForecast.Range("D7").Value = 312
$wsForecast.Range("H7").Value = 13.52
$wsForecast.Range("L7").Value = 0.92

# Row 8 (W16)
$wsForecast.Range("D8").Value = 290
$wsForecast.Range("H8").Value = 13.47
$wsForecast.Range("L8").Value = 0.87

# Row 9 (W17)
$wsForecast.Range("D9").Value = 266
$wsForecast.Range("H9").Value = 13.6

# Row 10 (W18)
$wsForecast.Range("D10").Value = 260
$wsForecast.Range("H10").Value = 12.89
$wsForecast.Range("L10").Value = 1

# Row 11 (W19)
$wsForecast.Range("D11").Value = 259
$wsForecast.Range("H11").Value = 11.93
$wsForecast.Range("L11").Value = 1.15

# Row 12 (W20)
$wsForecast.Range("D12").Value = 241
$wsForecast.Range("H12").Value = 11.75
$wsForecast.Range("L12").Value = 1.09

# Row 13 (W21)
$wsForecast.Range("D13").Value = 233
$wsForecast.Range("H13").Value = 11.12
$wsForecast.Range("L13").Value = 1

# Row 14 (W22)
$wsForecast.Range("D14").Value = 227
$wsForecast.Range("H14").Value = 10.39
$wsForecast.Range("J14").Value = "Normal"
$wsForecast.Range("L14").Value = 1.12

# Row 15 (W23)
$wsForecast.Range("D15").Value = 217
$wsForecast.Range("H15").Value = 9.82
$wsForecast.Range("I15").Value = "Low"
$wsForecast.Range("J15").Value = "Normal"
$wsForecast.Range("L15").Value = 1.07

# Row 16 (W24)
$wsForecast.Range("D16").Value = 232
$wsForecast.Range("H16").Value = 8.25
$wsForecast.Range("I16").Value = "Low"
$wsForecast.Range("J16").Value = "Normal"
$wsForecast.Range("L16").Value = 0.88

# Row 17 (W25)
$wsForecast.Range("D17").Value = 233
$wsForecast.Range("H17").Value = 7.22
$wsForecast.Range("I17").Value = "Low"
$wsForecast.Range("J17").Value = "Normal"
$wsForecast.Range("L17").Value = 1.08

# --- Summary sheet --------------------------------------------------------
# These "numbers" are stored as text on this sheet (same as the original
# file), so a leading apostrophe keeps them as text instead of Excel's
# normal auto-conversion to a numeric value.
$wsSummary.Range("B9").Value  = "'4736"
$wsSummary.Range("B10").Value = "'2834"
$wsSummary.Range("B11").Value = "'1627"
$wsSummary.Range("B12").Value = "'471"
$wsSummary.Range("B14").Value = "'217"
